# Update the lattice-multiplication exercise table: every cell keeps its
# 5-line layout (problem / top factor digits / rule / two partial-product
# rows) but the numbers themselves are regenerated.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# w:br line breaks inside a Range.Text assignment are represented by the
# vertical-tab control character (chr 11).
$vt = [char]11

function Set-ExerciseCell($table, $row, $col, $problem, $factors, $rule, $line4, $line5) {
    $table.Cell($row, $col).Range.Text = $problem + $vt + $factors + $vt + $rule + $vt + $line4 + $vt + $line5
}

Set-ExerciseCell $t 1 1 "43 x 45" "  4    5" "  ----" "4|    |" "3|    |"
Set-ExerciseCell $t 1 2 "89 x 22" "  2    2" "  ----" "8|    |" "9|    |"
Set-ExerciseCell $t 1 3 "98 x 39" "  3    9" "  ----" "9|    |" "8|    |"

Set-ExerciseCell $t 2 1 "92 x 25" "  2    5" "  ----" "9|    |" "2|    |"
Set-ExerciseCell $t 2 2 "22 x 95" "  9    5" "  ----" "2|    |" "2|    |"
Set-ExerciseCell $t 2 3 "97 x 68" "  6    8" "  ----" "9|    |" "7|    |"

Set-ExerciseCell $t 3 1 "24 x 74" "  7    4" "  ----" "2|    |" "4|    |"
Set-ExerciseCell $t 3 2 "63 x 51" "  5    1" "  ----" "6|    |" "3|    |"
Set-ExerciseCell $t 3 3 "33 x 76" "  7    6" "  ----" "3|    |" "3|    |"

Set-ExerciseCell $t 4 1 "46 x 24" "  2    4" "  ----" "4|    |" "6|    |"
Set-ExerciseCell $t 4 2 "97 x 47" "  4    7" "  ----" "9|    |" "7|    |"
Set-ExerciseCell $t 4 3 "30 x 76" "  7    6" "  ----" "3|    |" "0|    |"

Set-ExerciseCell $t 5 1 "64 x 93" "  9    3" "  ----" "6|    |" "4|    |"
Set-ExerciseCell $t 5 2 "44 x 27" "  2    7" "  ----" "4|    |" "4|    |"
Set-ExerciseCell $t 5 3 "70 x 90" "  9    0" "  ----" "7|    |" "0|    |"

Write-Host "Updated" $t.Rows.Count "x" $t.Columns.Count "exercise table"
